# Reorganized / exposed final API: populate the "NOTES" (column D) cells
# for the GIF header bytes (rows 2-11) and the trailer byte (row 54) of the
# byte-by-byte comparison table on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value  = "G"
$ws.Range("D3").Value  = "I"
$ws.Range("D4").Value  = "F"
$ws.Range("D5").Value  = 8
$ws.Range("D6").Value  = 9
$ws.Range("D7").Value  = "a"
$ws.Range("D8").Value  = "width:"
$ws.Range("D9").Value  = 3
$ws.Range("D10").Value = "height:"
$ws.Range("D11").Value = 3
$ws.Range("D54").Value = "trailer"

# Move the view/selection down to where the newly-filled-in rows are.
$ws.Range("J57").Select() | Out-Null
